# Refresh the cryptocurrency price/volume snapshot (Price = column D, Volume(1h) = column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.915.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.551.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.772.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.556.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.899.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.424.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +4.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.956"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.520"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.687.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0523"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("E51").Value = "  +1.68%  "
